# Update "Alimentadoras Estatales" (column E) values on the single
# worksheet of the workbook, per the updated source database ("Se subio
# la base de datos para el municipal"). Rows 48 and 71 are left untouched
# because their values did not change in the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = 140
    3 = 58
    4 = 52
    5 = 55
    6 = 4
    7 = 31
    8 = 71
    9 = 22
    10 = 0
    11 = 40
    12 = 0
    13 = 186
    14 = 29
    15 = 127
    16 = 45
    17 = 76
    18 = 0
    19 = 105
    20 = 33
    21 = 60
    22 = 0
    23 = 60
    24 = 0
    25 = 116
    26 = 152
    27 = 70
    28 = 128
    29 = 189
    30 = 20
    31 = 24
    32 = 134
    33 = 16
    34 = 44
    35 = 127
    36 = 69
    37 = 79
    38 = 168
    39 = 23
    40 = 99
    41 = 25
    42 = 0
    43 = 97
    44 = 156
    45 = 54
    46 = 31
    47 = 0
    49 = 106
    50 = 0
    51 = 80
    52 = 64
    53 = 126
    54 = 188
    55 = 0
    56 = 26
    57 = 13
    58 = 75
    59 = 42
    60 = 44
    61 = 105
    62 = 50
    63 = 103
    64 = 25
    65 = 5
    66 = 0
    67 = 0
    68 = 124
    69 = 0
    70 = 0
    72 = 6
    73 = 177
    74 = 3
    75 = 4
    76 = 5
    77 = 16
    78 = 19
    79 = 79
    80 = 85
    81 = 48
    82 = 97
    83 = 19
    84 = 61
    85 = 313
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}
